$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "From" value of rule R30 (row 10) changes from 18 to 1
$ws.Range("C10").Value = 1
